# Auto-generated: applies the cryptos.xlsx price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.561.51"
$ws.Range("E2").Value = "  +2.20%  "

# Row 3
$ws.Range("D3").Value = "2.686.95"
$ws.Range("E3").Value = "  -0.91%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.41%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.05"
$ws.Range("E5").Value = "  +1.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.94"
$ws.Range("E6").Value = "  -0.39%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("E8").Value = "  +1.36%  "

# Row 9
$ws.Range("D9").Value = "2.701.21"
$ws.Range("E9").Value = "  -0.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.45"
$ws.Range("E10").Value = "  +2.70%  "

# Row 11
$ws.Range("E11").Value = "  -0.90%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.339"
$ws.Range("E12").Value = "  -0.20%  "

# Row 13
$ws.Range("E13").Value = "  +2.61%  "

# Row 14
$ws.Range("D14").Value = "3.154.31"
$ws.Range("E14").Value = "  -0.75%  "

# Row 15
$ws.Range("D15").Value = "60.539.23"
$ws.Range("E15").Value = "  +2.07%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.33"
$ws.Range("E16").Value = "  +0.90%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("E17").Value = "  +0.19%  "

# Row 18
$ws.Range("D18").Value = "2.695.81"
$ws.Range("E18").Value = "  -0.35%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "349.82"
$ws.Range("E19").Value = "  -1.95%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.53"
$ws.Range("E20").Value = "  -0.87%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.58"
$ws.Range("E21").Value = "  +1.07%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.32"
$ws.Range("E22").Value = "  +1.35%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.58"
$ws.Range("E24").Value = "  +2.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.421"
$ws.Range("E25").Value = "  -0.52%  "

# Row 26
$ws.Range("E26").Value = "  +5.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("E27").Value = "  +0.35%  "

# Row 28
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.36"
$ws.Range("E28").Value = "  +1.48%  "

# Row 29
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0819"
$ws.Range("E29").Value = "  +0.30%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.91"
$ws.Range("E30").Value = "  +8.30%  "

# Row 31
$ws.Range("E31").Value = "  +0.12%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.27"
$ws.Range("E32").Value = "  +0.66%  "

# Row 33
$ws.Range("E33").Value = "  +0.64%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.85"
$ws.Range("E34").Value = "  -0.59%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.27"
$ws.Range("E35").Value = "  +5.37%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.25"
$ws.Range("E36").Value = "  +9.11%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.953"
$ws.Range("E37").Value = "  -2.65%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.883"
$ws.Range("E38").Value = "  +3.93%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.52"
$ws.Range("E39").Value = "  +7.07%  "

# Row 40
$ws.Range("E40").Value = "  +0.29%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.66"
$ws.Range("E41").Value = "  -2.35%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "283.40"
$ws.Range("E42").Value = "  +0.03%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.14"
$ws.Range("E43").Value = "  +1.19%  "

# Row 44
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0991"
$ws.Range("E44").Value = "  -0.08%  "

# Row 45
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.611"
$ws.Range("E45").Value = "  -1.54%  "

# Row 46
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.996"
$ws.Range("E46").Value = "  -0.24%  "

# Row 47
$ws.Range("D47").Value = "2.141.06"
$ws.Range("E47").Value = "  +6.12%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.90"
$ws.Range("E48").Value = "  +3.28%  "

# Row 49
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0539"
$ws.Range("E49").Value = "  +0.87%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0235"
$ws.Range("E50").Value = "  +1.12%  "

# Row 51
$ws.Range("E51").Value = "  +1.65%  "
